$d = $word.ActiveDocument

# Locate the paragraph containing the "Ver no Jupiter ..." text, which marks
# the start (together with the blank paragraph right before it) of the block
# that needs to be removed. The block to delete runs from the blank paragraph
# just after the SAMUELSON bibliography entry through the "(c) 2020 ..."
# footer paragraph (inclusive), leaving the following blank paragraph and the
# page-break paragraph untouched.

$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$jupiterIndex = 0
$copyrightIndex = 0

for ($i = 1; $i -le $count; $i++) {
    $text = $paragraphs.Item($i).Range.Text
    if ($text -like "Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    if ($text -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -ge $jupiterIndex) {
    # Include the blank paragraph immediately preceding the "Ver no Jupiter"
    # paragraph so it is removed too.
    $startIndex = $jupiterIndex - 1
    if ($startIndex -lt 1) {
        $startIndex = $jupiterIndex
    }

    $startRange = $paragraphs.Item($startIndex).Range
    $endRange = $paragraphs.Item($copyrightIndex).Range

    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}
